# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rebuilds the detail table (rows 16-40) on "Hoja1":
#   - 16 new rows for DIANA CAROLINA GONZALEZ OQUINES (CC 1128063396) covering
#     periods 2008-2111, each with Valor Mora 39227 and Salario Basico 980657.
#   - The pre-existing BETTY MARCELA PANTOJA AGAMEZ (CC 1143327342) row for
#     period 2112 (Valor Mora 21804 / Salario Basico 1160000) is kept, now
#     relocated to row 32.
#   - 8 more DIANA rows for periods 2112, 2201-2207 (same amounts as above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2008", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2009", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2010", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2011", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2012", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2101", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2102", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2103", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2104", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2105", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2106", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2107", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2108", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2109", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2110", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2111", 39227, 980657)
    ,@("CC", "1143327342", "BETTY MARCELA PANTOJA AGAMEZ", "2112", 21804, 1160000)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2112", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2201", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2202", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2203", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2204", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2205", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2206", 39227, 980657)
    ,@("CC", "1128063396", "DIANA CAROLINA GONZALEZ OQUINES", "2207", 39227, 980657)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 2).Value2 = $rec[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value2 = $rec[1]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value2 = $rec[2]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value2 = $rec[3]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value2 = $rec[4]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value2 = $rec[5]   # G - Salario Basico
}
